# Modified the Test Methods for Customer portal
# Delete the "Verify Withdraw Token with New and Existing Debit Cards" test row (row 10),
# and update the iteration counts in row 9 (Stop/Start Iteration) from 9 to 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 Start/Stop Iteration values from 9 to 14 (keep as quote-prefixed text
# to preserve the existing "quotePrefix" cell style used for these numeric-looking text values)
$ws.Range("E9").Value = "'14"
$ws.Range("F9").Value = "'14"

# Delete entire row 10 (shifts nothing below it, it's the last row)
$ws.Rows("10:10").Delete()

# Update the view state: selected cell / top-left cell
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B2:B9").Select()
